$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false
$ws = $wb.Worksheets.Item("Estado Ruta")
$ws.Delete()
Write-Output "Deleted"
Write-Output $wb.Worksheets.Count
